# Rename the worksheet from "ValidLogin" to "LoginData"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidLogin")
$ws.Name = "LoginData"
